$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Correct the objective function value (named range "b" = Sheet1!$B$4) ---
$ws.Range("B4").Value = 8

# --- Zero out the "x" matrix (named range "x" = Sheet1!$C$20:$I$26) ---
for ($row = 20; $row -le 26; $row++) {
    for ($col = 3; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}

# --- Update the "t" matrix (named range "t" = Sheet1!$C$28:$I$34) with corrected values ---
$tValues = @{
    28 = @(0,  4,  8,  12, 14, 16, 18)
    29 = @(10, 14, 16, 19, 21, 24, 26)
    30 = @(20, 22, 26, 28, 30, 34, 33)
    31 = @(30, 32, 36, 38, 40, 42, 42)
    32 = @(39, 41, 45, 47, 51, 52, 53)
    33 = @(48, 52, 54, 58, 60, 62, 63)
    34 = @(57, 61, 64, 68, 70, 72, 74)
}

foreach ($row in $tValues.Keys) {
    $rowValues = $tValues[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = 3 + $i
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}
